# Apply the commit "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# 1. Update VALOR MORA total (E11): 180241 -> 66361
# 2. Update Cant. Trabajadores (C13): 4 -> 2
# 3. Update Cant. Periodos (F13): 3 -> 2
# 4. Remove the two obsolete worker rows (18 and 19: CARLOS ALFREDO VIANA
#    MONTEROSA / 1120740842 / 2508 and NELSON ARTURO TORRECILLA MOLINA /
#    19874875 / 2508), shifting the footer rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = 66361
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

$ws.Rows("18:19").Delete()

# Excel auto-applied center alignment to the "Periodo Mora" column for the
# two surviving worker rows once the table shrank to two rows.
$ws.Range("E16:E17").HorizontalAlignment = -4108
